# Daily attendance processing - 2026-01-28 05:20:50
# Reorders the comma-separated "Recorded By" list in column G so that any
# entry that is exactly "System" is moved to the end of the list, while the
# remaining entries are sorted alphabetically (case-insensitive).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    if ($val -notmatch ',') { continue }

    $parts = $val -split ',' | ForEach-Object { $_.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.Equals('System')) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $sortedOthers = $otherParts | Sort-Object { $_.ToLower() }

    $newParts = @()
    $newParts += $sortedOthers
    $newParts += $systemParts

    $newVal = [string]::Join(', ', $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
